$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Invoice Number"
$ws.Range("B1").Value = "Vendor"
$ws.Range("C1").Value = "Vendor Code"
$ws.Range("D1").Value = "Service"
$ws.Range("E1").Value = "Date"
$ws.Range("F1").Value = "Total Amount"
$ws.Range("G1").Value = "Description"

# Data row
$ws.Range("A2").Value = "INV-2025-001"
$ws.Range("B2").Value = "NovaTech Solutions"
$ws.Range("C2").Value = "VEND001"
$ws.Range("D2").Value = "DELIVERY"

# Keep the date cell as literal text (not an Excel date serial number)
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "01/15/2024"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = 250
$ws.Range("G2").Value = "Professional Services"

# F1/G1 are new header cells outside the old used range -- they don't inherit the
# bold/centered/bordered header style automatically, so copy it explicitly from A1.
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths (closest values reachable through the host's pixel-snapped column sizing)
$ws.Columns("A").ColumnWidth = 19
$ws.Range("B1:D1").EntireColumn.ColumnWidth = 18
$ws.Columns("E").ColumnWidth = 14.3
$ws.Columns("F").ColumnWidth = 11.65
$ws.Columns("G").ColumnWidth = 17.3
$ws.Columns("H").ColumnWidth = 12.3

# Restore the selection to match the saved view
[void]$ws.Range("G13").Select()
